$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3571.611
$ws.Range("I51").Value = 3842.7144
$ws.Range("J51").Value = 3399.0908
$ws.Range("K51").Value = 3842.7144
$ws.Range("L51").Value = 3399.0908
$ws.Range("M51").Value = -3358.7144
$ws.Range("N51").Value = -4367.0908
$ws.Range("H64").Value = 3199.1667
$ws.Range("I64").Value = 2750
$ws.Range("K64").Value = 2750
$ws.Range("M64").Value = -2502
$ws.Range("H67").Value = 3199.1667
$ws.Range("I67").Value = 2750
$ws.Range("K67").Value = 2750
$ws.Range("M67").Value = -1892
$ws.Range("H74").Value = 3930.7693
$ws.Range("I74").Value = 3816.6667
$ws.Range("J74").Value = 3965
$ws.Range("K74").Value = 3816.6667
$ws.Range("L74").Value = 3965
$ws.Range("M74").Value = -2880.6667
$ws.Range("N74").Value = -5837
$ws.Range("H77").Value = 3930.7693
$ws.Range("I77").Value = 3816.6667
$ws.Range("J77").Value = 3965
$ws.Range("K77").Value = 19083.3335
$ws.Range("L77").Value = 19825
$ws.Range("M77").Value = -14403.3335
$ws.Range("N77").Value = -29185
$ws.Range("H106").Value = 20041720
$ws.Range("I106").Value = 48999.047
$ws.Range("J106").Value = 125003500
$ws.Range("K106").Value = 48999.047
$ws.Range("L106").Value = 125003500
$ws.Range("M106").Value = -48368.047
$ws.Range("N106").Value = -125004762
$ws.Range("H137").Value = 1410.32
$ws.Range("J137").Value = 2360
$ws.Range("L137").Value = 7080
$ws.Range("N137").Value = -12180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17791.688
$ws.Range("I32").Value = 18963.803
$ws.Range("K32").Value = 18963.803
$ws.Range("M32").Value = -18676.803
$ws.Range("H61").Value = 1493.826
$ws.Range("I61").Value = 1154.7059
$ws.Range("J61").Value = 2454.6667
$ws.Range("K61").Value = 1154.7059
$ws.Range("L61").Value = 2454.6667
$ws.Range("M61").Value = -942.7058999999999
$ws.Range("N61").Value = -2878.6667
$ws.Range("H96").Value = 25947.834
$ws.Range("J96").Value = 25947.834
$ws.Range("L96").Value = 25947.834
$ws.Range("N96").Value = -31439.834
$ws.Range("H132").Value = 4335.959
$ws.Range("I132").Value = 4653.324
$ws.Range("K132").Value = 13959.972
$ws.Range("M132").Value = -11429.972
$ws.Range("H136").Value = 1493.826
$ws.Range("I136").Value = 1154.7059
$ws.Range("J136").Value = 2454.6667
$ws.Range("K136").Value = 3464.1177
$ws.Range("L136").Value = 7364.000100000001
$ws.Range("M136").Value = -914.1176999999998
$ws.Range("N136").Value = -12464.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1863.8572
$ws.Range("I99").Value = 1674.5
$ws.Range("K99").Value = 1674.5
$ws.Range("M99").Value = -176.5
$ws.Range("H107").Value = 1618.8334
$ws.Range("I107").Value = 1754.0834
$ws.Range("K107").Value = 1754.0834
$ws.Range("M107").Value = 165.9166
$ws.Range("H134").Value = 1886.6731
$ws.Range("I134").Value = 1665.3529
$ws.Range("J134").Value = 2304.7222
$ws.Range("K134").Value = 4996.0587
$ws.Range("L134").Value = 6914.1666
$ws.Range("M134").Value = -2461.0587
$ws.Range("N134").Value = -11984.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4905361
$ws.Range("I31").Value = 2662.926
$ws.Range("J31").Value = 23815768
$ws.Range("K31").Value = 2662.926
$ws.Range("L31").Value = 23815768
$ws.Range("M31").Value = -2367.926
$ws.Range("N31").Value = -23816358
$ws.Range("H34").Value = 4905361
$ws.Range("I34").Value = 2662.926
$ws.Range("J34").Value = 23815768
$ws.Range("K34").Value = 2662.926
$ws.Range("L34").Value = 23815768
$ws.Range("M34").Value = -2460.926
$ws.Range("N34").Value = -23816172
$ws.Range("H86").Value = 90910920
$ws.Range("I86").Value = 166668270
$ws.Range("J86").Value = 2093
$ws.Range("K86").Value = 166668270
$ws.Range("L86").Value = 2093
$ws.Range("M86").Value = -166667147
$ws.Range("N86").Value = -4339
$ws.Range("H89").Value = 90910920
$ws.Range("I89").Value = 166668270
$ws.Range("J89").Value = 2093
$ws.Range("K89").Value = 833341350
$ws.Range("L89").Value = 10465
$ws.Range("M89").Value = -833335734
$ws.Range("N89").Value = -21697
$ws.Range("H99").Value = 2602.375
$ws.Range("I99").Value = 2499.1428
$ws.Range("J99").Value = 2682.6667
$ws.Range("K99").Value = 2499.1428
$ws.Range("L99").Value = 2682.6667
$ws.Range("M99").Value = -1001.1428
$ws.Range("N99").Value = -5678.6667
$ws.Range("H126").Value = 2602.375
$ws.Range("I126").Value = 2499.1428
$ws.Range("J126").Value = 2682.6667
$ws.Range("K126").Value = 7497.428400000001
$ws.Range("L126").Value = 8048.000100000001
$ws.Range("M126").Value = -5027.428400000001
$ws.Range("N126").Value = -12988.0001
$ws.Range("H134").Value = 1228.6316
$ws.Range("I134").Value = 1185.8667
$ws.Range("K134").Value = 3557.6001
$ws.Range("M134").Value = -1022.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1297.5
$ws.Range("I5").Value = 264.63635
$ws.Range("J5").Value = 1718.2963
$ws.Range("K5").Value = 793.90905
$ws.Range("L5").Value = 5154.8889
$ws.Range("M5").Value = -681.90905
$ws.Range("N5").Value = -5378.8889
$ws.Range("H131").Value = 814.96
$ws.Range("I131").Value = 625
$ws.Range("J131").Value = 827.0851
$ws.Range("K131").Value = 1875
$ws.Range("L131").Value = 2481.2553
$ws.Range("M131").Value = 3165
$ws.Range("N131").Value = -12561.2553
$ws.Range("H135").Value = 1297.5
$ws.Range("I135").Value = 264.63635
$ws.Range("J135").Value = 1718.2963
$ws.Range("K135").Value = 2381.72715
$ws.Range("L135").Value = 15464.6667
$ws.Range("M135").Value = 153.2728500000003
$ws.Range("N135").Value = -20534.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 250002260
$ws.Range("J92").Value = 250002260
$ws.Range("L92").Value = 250002260
$ws.Range("N92").Value = -250006004
$ws.Range("H97").Value = 816.55554
$ws.Range("I97").Value = 693.9545000000001
$ws.Range("J97").Value = 1356
$ws.Range("K97").Value = 693.9545000000001
$ws.Range("L97").Value = 1356
$ws.Range("M97").Value = -197.9545000000001
$ws.Range("N97").Value = -2348
$ws.Range("H113").Value = 17857794
$ws.Range("I113").Value = 62500240
$ws.Range("J113").Value = 816.2
$ws.Range("K113").Value = 62500240
$ws.Range("L113").Value = 816.2
$ws.Range("M113").Value = -62498070
$ws.Range("N113").Value = -5156.2
$ws.Range("H132").Value = 24930.023
$ws.Range("I132").Value = 33045.47
$ws.Range("J132").Value = 3288.8333
$ws.Range("K132").Value = 99136.41
$ws.Range("L132").Value = 9866.499899999999
$ws.Range("M132").Value = -96606.41
$ws.Range("N132").Value = -14926.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2211.5557
$ws.Range("I93").Value = 2250
$ws.Range("J93").Value = 2200.5715
$ws.Range("K93").Value = 2250
$ws.Range("L93").Value = 2200.5715
$ws.Range("M93").Value = -1002
$ws.Range("N93").Value = -4696.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3599.2104
$ws.Range("I81").Value = 3153.2222
$ws.Range("J81").Value = 4000.6
$ws.Range("K81").Value = 6306.4444
$ws.Range("L81").Value = 8001.2
$ws.Range("M81").Value = -5245.4444
$ws.Range("N81").Value = -10123.2
$ws.Range("H84").Value = 3599.2104
$ws.Range("I84").Value = 3153.2222
$ws.Range("J84").Value = 4000.6
$ws.Range("K84").Value = 31532.222
$ws.Range("L84").Value = 40006
$ws.Range("M84").Value = -26228.222
$ws.Range("N84").Value = -50614
$ws.Range("H92").Value = 17545
$ws.Range("J92").Value = 17545
$ws.Range("L92").Value = 17545
$ws.Range("N92").Value = -22537
$ws.Range("H96").Value = 83333700
$ws.Range("I96").Value = 83333700
$ws.Range("K96").Value = 83333700
$ws.Range("M96").Value = -83332327
